$d = $word.ActiveDocument

# The "Professional Summary" opening sentence changes from:
#   "Desenvolvedor Full Stack atuando em desenvolvimento de software desde 2020"
# to:
#   "Atuando no desenvolvimento de softwares desde 2020"
# (the rest of the paragraph - starting at ". Especializado ..." - is untouched)

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Desenvolvedor Full Stack atuando em desenvolvimento de software desde 2020",
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "Atuando no desenvolvimento de softwares desde 2020", # ReplaceWith
    2        # Replace (wdReplaceAll)
)
